$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: write new values as text-forcing formulas (="...") so Excel
# does not auto-coerce numeric-looking strings (e.g. "214.94") into
# real numbers -- these cells must stay text, matching the source data.
$ws.Range("D2").Formula = '="27.024.20"'
$ws.Range("E2").Formula = '="  +2.85%  "'
$ws.Range("D3").Formula = '="1.649.99"'
$ws.Range("E3").Formula = '="  +3.51%  "'
$ws.Range("D5").Formula = '="214.94"'
$ws.Range("E5").Formula = '="  +1.68%  "'
$ws.Range("E6").Formula = '="  +1.56%  "'
$ws.Range("E7").Formula = '="  -0.02%  "'
$ws.Range("E8").Formula = '="  +1.88%  "'
$ws.Range("E9").Formula = '="  +1.59%  "'
$ws.Range("E10").Formula = '="  +3.94%  "'
$ws.Range("D11").Formula = '="0.0866"'
$ws.Range("E11").Formula = '="  +1.21%  "'
$ws.Range("D12").Formula = '="1.882.82"'
$ws.Range("E12").Formula = '="  +3.53%  "'
$ws.Range("D13").Formula = '="1.650.74"'
$ws.Range("E13").Formula = '="  +3.51%  "'
$ws.Range("E14").Formula = '="  +2.46%  "'
$ws.Range("E15").Formula = '="  +3.11%  "'
$ws.Range("D16").Formula = '="65.08"'
$ws.Range("E16").Formula = '="  +2.61%  "'
$ws.Range("D17").Formula = '="27.001.49"'
$ws.Range("E17").Formula = '="  +2.79%  "'
$ws.Range("D18").Formula = '="238.24"'
$ws.Range("E18").Formula = '="  +4.22%  "'
$ws.Range("D19").Formula = '="7.85"'
$ws.Range("E19").Formula = '="  +2.34%  "'
$ws.Range("E20").Formula = '="  +1.30%  "'
$ws.Range("E21").Formula = '="  +0.05%  "'
$ws.Range("E22").Formula = '="  +4.41%  "'
$ws.Range("E23").Formula = '="  +4.45%  "'
$ws.Range("D24").Formula = '="9.22"'
$ws.Range("E24").Formula = '="  +3.53%  "'
$ws.Range("D25").Formula = '="145.43"'
$ws.Range("E25").Formula = '="  -0.52%  "'
$ws.Range("D28").Formula = '="0.114"'
$ws.Range("E28").Formula = '="  +1.55%  "'
$ws.Range("E29").Formula = '="  +2.72%  "'
$ws.Range("E30").Formula = '="  +0.79%  "'
$ws.Range("E31").Formula = '="  +1.85%  "'
$ws.Range("E32").Formula = '="  +3.30%  "'
$ws.Range("D33").Formula = '="1.509.16"'
$ws.Range("E33").Formula = '="  +2.41%  "'
$ws.Range("E34").Formula = '="  +5.02%  "'
$ws.Range("D36").Formula = '="2.41"'
$ws.Range("E36").Formula = '="  -0.16%  "'
$ws.Range("D37").Formula = '="0.577"'
$ws.Range("E37").Formula = '="  +1.66%  "'
$ws.Range("D38").Formula = '="0.886"'
$ws.Range("E38").Formula = '="  +8.71%  "'
$ws.Range("D39").Formula = '="0.0168"'
$ws.Range("E39").Formula = '="  +2.57%  "'
$ws.Range("E40").Formula = '="  +3.46%  "'
$ws.Range("E41").Formula = '="  -0.01%  "'
$ws.Range("E42").Formula = '="  +4.31%  "'
$ws.Range("E43").Formula = '="  +9.56%  "'
$ws.Range("D44").Formula = '="1.790.01"'
$ws.Range("E44").Formula = '="  +3.32%  "'
$ws.Range("D45").Formula = '="0.775"'
$ws.Range("E45").Formula = '="  +2.79%  "'
$ws.Range("E46").Formula = '="  -1.35%  "'
$ws.Range("D47").Formula = '="89.45"'
$ws.Range("E47").Formula = '="  +0.76%  "'
$ws.Range("E48").Formula = '="  -0.08%  "'
$ws.Range("E49").Formula = '="  +3.16%  "'
$ws.Range("D50").Formula = '="0.0506"'
$ws.Range("E50").Formula = '="  +1.28%  "'
$ws.Range("D51").Formula = '="0.0974"'
$ws.Range("E51").Formula = '="  +2.16%  "'
# Step 2: convert the helper formulas back into plain static values
# (Copy + PasteSpecial values-only) so no formula or extra style survives.
$rng = $ws.Range("D2:E51")
$rng.Copy()
$rng.PasteSpecial(-4163)

Write-Host "Updated cryptos list"
